$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("D8:D17").Value = 0.0001
